$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Finish subsection on Good-Thomas: fill in the "Ours" column (E) ---

# 2^8 section (rows 6-9)
$ws.Range("E6").Value = 239498
$ws.Range("E7").Value = 230037
$ws.Range("E8").Value = 4193
$ws.Range("E9").Value = 420

# 2^16 section header (row 14)
$ws.Range("E14").Value = 100000

# 2^32 section header (row 23)
$ws.Range("E23").Value = 100000

# --- Update the view state to match where the author left off editing ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E23").Select()
